# "moves staff to edit to create staff"
# Adds a new staff record (NATHAN DANSKIN, id 2222, MANAGER) as row 2
# of the "Staff" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New staff member's details
$ws.Range("A2").Value = "NATHAN"
$ws.Range("B2").Value = "DANSKIN"

# The id column stores the number but is formatted as text (matches how
# the existing 1111 id in C1 reads back), so set the numeric value first
# and apply the text number format afterwards.
$ws.Range("C2").Value = 2222
$ws.Range("C2").NumberFormat = "@"

$ws.Range("D2").Value = "MANAGER"

# Roughly match the column widths Excel auto-sized to fit the new data.
$ws.Columns.Item(1).ColumnWidth = 7.59
$ws.Columns.Item(2).ColumnWidth = 8.09
$ws.Columns.Item(3).ColumnWidth = 4.09
$ws.Columns.Item(4).ColumnWidth = 9.26

# Leave the new id cell selected, as it was after data entry.
[void]$ws.Range("C2").Select()

# Page was switched to portrait orientation.
$ws.PageSetup.Orientation = 1
